$wb = $excel.ActiveWorkbook

$wsIEEE = $wb.Worksheets.Item("NetworkLine_IEEE")
$wsIEEE.Range("D8").Select()

$wsNL = $wb.Worksheets.Item("NetworkLine")
$wsNL.Range("D15").Select()
